$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 3070.125
$ws.Cells.Item(2, 9).Value = 316
$ws.Cells.Item(2, 11).Value = 316
$ws.Cells.Item(2, 13).Value = -203
$ws.Cells.Item(8, 8).Value = 558.1667
$ws.Cells.Item(8, 9).Value = 110.8
$ws.Cells.Item(8, 11).Value = 332.4
$ws.Cells.Item(8, 13).Value = -193.4
$ws.Cells.Item(32, 8).Value = 3602.625
$ws.Cells.Item(32, 9).Value = 3464.5
$ws.Cells.Item(32, 10).Value = 3648.6667
$ws.Cells.Item(32, 11).Value = 3464.5
$ws.Cells.Item(32, 12).Value = 3648.6667
$ws.Cells.Item(32, 13).Value = -3138.5
$ws.Cells.Item(32, 14).Value = -4300.6667
$ws.Cells.Item(38, 8).Value = 3100.7693
$ws.Cells.Item(38, 9).Value = 2534.1667
$ws.Cells.Item(38, 11).Value = 7602.500100000001
$ws.Cells.Item(38, 13).Value = -7230.500100000001
$ws.Cells.Item(58, 8).Value = 430
$ws.Cells.Item(58, 9).Value = 430
$ws.Cells.Item(58, 11).Value = 1290
$ws.Cells.Item(58, 13).Value = -1140
$ws.Cells.Item(107, 8).Value = 320.33334
$ws.Cells.Item(107, 9).Value = 320.33334
$ws.Cells.Item(107, 11).Value = 320.33334
$ws.Cells.Item(107, 13).Value = 1599.66666
$ws.Cells.Item(125, 8).Value = 482
$ws.Cells.Item(125, 9).Value = 482
$ws.Cells.Item(125, 10).Value = 0
$ws.Cells.Item(125, 11).Value = 4338
$ws.Cells.Item(125, 12).ClearContents()
$ws.Cells.Item(125, 13).Value = -1878
$ws.Cells.Item(125, 14).Value = 0
$ws.Cells.Item(135, 8).Value = 1859.6111
$ws.Cells.Item(135, 9).Value = 2071.7144
$ws.Cells.Item(135, 11).Value = 18645.4296
$ws.Cells.Item(135, 13).Value = -16110.4296
$ws.Cells.Item(137, 8).Value = 1269.1333
$ws.Cells.Item(137, 9).Value = 1050
$ws.Cells.Item(137, 10).Value = 1597.8334
$ws.Cells.Item(137, 11).Value = 3150
$ws.Cells.Item(137, 12).Value = 4793.5002
$ws.Cells.Item(137, 13).Value = -600
$ws.Cells.Item(137, 14).Value = -9893.5002
$ws.Cells.Item(138, 8).Value = 5597.3184
$ws.Cells.Item(138, 9).Value = 5198
$ws.Cells.Item(138, 11).Value = 15594
$ws.Cells.Item(138, 13).Value = -10454
$ws.Cells.Item(141, 8).Value = 3499.6
$ws.Cells.Item(141, 9).Value = 2999.3333
$ws.Cells.Item(141, 11).Value = 8997.999899999999
$ws.Cells.Item(141, 13).Value = -3817.999899999999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 3503533.8
$ws.Cells.Item(32, 9).Value = 3337035.5
$ws.Cells.Item(32, 11).Value = 3337035.5
$ws.Cells.Item(32, 13).Value = -3336748.5
$ws.Cells.Item(43, 8).Value = 19997.5
$ws.Cells.Item(43, 10).Value = 19997.5
$ws.Cells.Item(43, 12).Value = 19997.5
$ws.Cells.Item(43, 14).Value = -20623.5
$ws.Cells.Item(61, 8).Value = 2186.6
$ws.Cells.Item(61, 10).Value = 1777
$ws.Cells.Item(61, 12).Value = 1777
$ws.Cells.Item(61, 14).Value = -2201
$ws.Cells.Item(63, 8).Value = 12513.833
$ws.Cells.Item(63, 9).Value = 12513.833
$ws.Cells.Item(63, 11).Value = 12513.833
$ws.Cells.Item(63, 13).Value = -11827.833
$ws.Cells.Item(66, 8).Value = 12513.833
$ws.Cells.Item(66, 9).Value = 12513.833
$ws.Cells.Item(66, 11).Value = 62569.165
$ws.Cells.Item(66, 13).Value = -59137.165
$ws.Cells.Item(74, 8).Value = 1799
$ws.Cells.Item(74, 9).Value = 1799
$ws.Cells.Item(74, 11).Value = 1799
$ws.Cells.Item(74, 13).Value = -925
$ws.Cells.Item(77, 8).Value = 1799
$ws.Cells.Item(77, 9).Value = 1799
$ws.Cells.Item(77, 11).Value = 8995
$ws.Cells.Item(77, 13).Value = -4627
$ws.Cells.Item(132, 8).Value = 1055.2858
$ws.Cells.Item(132, 9).Value = 1022
$ws.Cells.Item(132, 10).Value = 1099.6666
$ws.Cells.Item(132, 11).Value = 3066
$ws.Cells.Item(132, 12).Value = 3298.9998
$ws.Cells.Item(132, 13).Value = -536
$ws.Cells.Item(132, 14).Value = -8358.9998
$ws.Cells.Item(133, 8).Value = 147000
$ws.Cells.Item(133, 10).Value = 150000
$ws.Cells.Item(133, 12).Value = 150000
$ws.Cells.Item(133, 14).Value = -155060
$ws.Cells.Item(136, 8).Value = 2186.6
$ws.Cells.Item(136, 10).Value = 1777
$ws.Cells.Item(136, 12).Value = 5331
$ws.Cells.Item(136, 14).Value = -10431

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 3307.1538
$ws.Cells.Item(134, 9).Value = 3307.1538
$ws.Cells.Item(134, 11).Value = 9921.4614
$ws.Cells.Item(134, 13).Value = -7386.4614

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2211.6206
$ws.Cells.Item(31, 10).Value = 2516.889
$ws.Cells.Item(31, 12).Value = 2516.889
$ws.Cells.Item(31, 14).Value = -3106.889
$ws.Cells.Item(34, 8).Value = 2211.6206
$ws.Cells.Item(34, 10).Value = 2516.889
$ws.Cells.Item(34, 12).Value = 2516.889
$ws.Cells.Item(34, 14).Value = -2920.889
$ws.Cells.Item(58, 8).Value = 3913.2
$ws.Cells.Item(58, 9).Value = 2947.25
$ws.Cells.Item(58, 11).Value = 2947.25
$ws.Cells.Item(58, 13).Value = -2744.25
$ws.Cells.Item(64, 8).Value = 0
$ws.Cells.Item(64, 10).Value = 0
$ws.Cells.Item(64, 12).ClearContents()
$ws.Cells.Item(64, 14).Value = 0
$ws.Cells.Item(67, 8).Value = 0
$ws.Cells.Item(67, 10).Value = 0
$ws.Cells.Item(67, 12).ClearContents()
$ws.Cells.Item(67, 14).Value = 0
$ws.Cells.Item(136, 8).Value = 3913.2
$ws.Cells.Item(136, 9).Value = 2947.25
$ws.Cells.Item(136, 11).Value = 8841.75
$ws.Cells.Item(136, 13).Value = -6291.75

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 89.5
$ws.Cells.Item(12, 10).Value = 89.5
$ws.Cells.Item(12, 12).Value = 268.5
$ws.Cells.Item(12, 14).Value = -614.5
$ws.Cells.Item(38, 8).Value = 2349.125
$ws.Cells.Item(38, 9).Value = 2980.5
$ws.Cells.Item(38, 10).Value = 455
$ws.Cells.Item(38, 11).Value = 8941.5
$ws.Cells.Item(38, 12).Value = 1365
$ws.Cells.Item(38, 13).Value = -8594.5
$ws.Cells.Item(38, 14).Value = -2059
$ws.Cells.Item(121, 8).Value = 3257.1765
$ws.Cells.Item(121, 9).Value = 917.5
$ws.Cells.Item(121, 10).Value = 4533.364
$ws.Cells.Item(121, 11).Value = 2752.5
$ws.Cells.Item(121, 12).Value = 13600.092
$ws.Cells.Item(121, 13).Value = -1442.5
$ws.Cells.Item(121, 14).Value = -16220.092
$ws.Cells.Item(131, 8).Value = 2110.15
$ws.Cells.Item(131, 10).Value = 2699.077
$ws.Cells.Item(131, 12).Value = 8097.231000000001
$ws.Cells.Item(131, 14).Value = -18177.231

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 243.57143
$ws.Cells.Item(2, 9).Value = 225.83333
$ws.Cells.Item(2, 11).Value = 225.83333
$ws.Cells.Item(2, 13).Value = -112.83333
$ws.Cells.Item(15, 8).Value = 34994.668
$ws.Cells.Item(15, 10).Value = 34994.668
$ws.Cells.Item(15, 12).Value = 34994.668
$ws.Cells.Item(15, 14).Value = -35570.668
$ws.Cells.Item(81, 8).Value = 34994.668
$ws.Cells.Item(81, 10).Value = 34994.668
$ws.Cells.Item(81, 12).Value = 34994.668
$ws.Cells.Item(81, 14).Value = -36990.668
$ws.Cells.Item(84, 8).Value = 34994.668
$ws.Cells.Item(84, 10).Value = 34994.668
$ws.Cells.Item(84, 12).Value = 104984.004
$ws.Cells.Item(84, 14).Value = -114968.004
$ws.Cells.Item(102, 8).Value = 1891.9286
$ws.Cells.Item(102, 9).Value = 1927.6
$ws.Cells.Item(102, 10).Value = 1802.75
$ws.Cells.Item(102, 11).Value = 1927.6
$ws.Cells.Item(102, 12).Value = 1802.75
$ws.Cells.Item(102, 13).Value = -305.5999999999999
$ws.Cells.Item(102, 14).Value = -5046.75
$ws.Cells.Item(122, 8).Value = 2628.1428
$ws.Cells.Item(122, 9).Value = 2758.2104
$ws.Cells.Item(122, 10).Value = 1392.5
$ws.Cells.Item(122, 11).Value = 8274.6312
$ws.Cells.Item(122, 12).Value = 4177.5
$ws.Cells.Item(122, 13).Value = -5824.6312
$ws.Cells.Item(122, 14).Value = -9077.5
$ws.Cells.Item(126, 8).Value = 8956.666999999999
$ws.Cells.Item(126, 10).Value = 8710
$ws.Cells.Item(126, 12).Value = 26130
$ws.Cells.Item(126, 14).Value = -31070
$ws.Cells.Item(132, 8).Value = 4899.9
$ws.Cells.Item(132, 9).Value = 5624.875
$ws.Cells.Item(132, 11).Value = 16874.625
$ws.Cells.Item(132, 13).Value = -14344.625

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 2823
$ws.Cells.Item(22, 9).Value = 2229.8572
$ws.Cells.Item(22, 10).Value = 4899
$ws.Cells.Item(22, 11).Value = 2229.8572
$ws.Cells.Item(22, 12).Value = 4899
$ws.Cells.Item(22, 13).Value = -1934.8572
$ws.Cells.Item(22, 14).Value = -5489
$ws.Cells.Item(27, 8).Value = 2823
$ws.Cells.Item(27, 9).Value = 2229.8572
$ws.Cells.Item(27, 10).Value = 4899
$ws.Cells.Item(27, 11).Value = 2229.8572
$ws.Cells.Item(27, 12).Value = 4899
$ws.Cells.Item(27, 13).Value = -2122.8572
$ws.Cells.Item(27, 14).Value = -5113
$ws.Cells.Item(46, 8).Value = 2188.182
$ws.Cells.Item(46, 9).Value = 1418.25
$ws.Cells.Item(46, 10).Value = 4241.3335
$ws.Cells.Item(46, 11).Value = 1418.25
$ws.Cells.Item(46, 12).Value = 4241.3335
$ws.Cells.Item(46, 13).Value = -1230.25
$ws.Cells.Item(46, 14).Value = -4617.3335
$ws.Cells.Item(132, 8).Value = 4149
$ws.Cells.Item(132, 10).Value = 3799
$ws.Cells.Item(132, 12).Value = 11397
$ws.Cells.Item(132, 14).Value = -16457

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(15, 8).Value = 39957.25
$ws.Cells.Item(15, 10).Value = 39957.25
$ws.Cells.Item(15, 12).Value = 39957.25
$ws.Cells.Item(15, 14).Value = -40533.25
$ws.Cells.Item(81, 8).Value = 6667001
$ws.Cells.Item(81, 9).Value = 0
$ws.Cells.Item(81, 11).Value = 0
$ws.Cells.Item(81, 13).ClearContents()
$ws.Cells.Item(84, 8).Value = 6667001
$ws.Cells.Item(84, 9).Value = 0
$ws.Cells.Item(84, 11).Value = 0
$ws.Cells.Item(84, 13).ClearContents()
$ws.Cells.Item(100, 8).Value = 11112617
$ws.Cells.Item(100, 9).Value = 12501611
$ws.Cells.Item(100, 10).Value = 670
$ws.Cells.Item(100, 11).Value = 25003222
$ws.Cells.Item(100, 12).Value = 1340
$ws.Cells.Item(100, 13).Value = -25002681
$ws.Cells.Item(100, 14).Value = -2422
$ws.Cells.Item(136, 8).Value = 2324.2
$ws.Cells.Item(136, 9).Value = 2325
$ws.Cells.Item(136, 10).Value = 2323
$ws.Cells.Item(136, 11).Value = 6975
$ws.Cells.Item(136, 12).Value = 6969
$ws.Cells.Item(136, 13).Value = -4425
$ws.Cells.Item(136, 14).Value = -12069
